# "Changed and fixed John Murungi"
#
# The DECEMBER 21 rent statement had three tenants (rows 7, 12, 13) whose
# "PAID" (column H) cells were left blank even though they had in fact
# paid in full. Fill in the amounts that were actually paid so the
# dependent "BAL" (column I) formulas - and everything that rolls up from
# them (row 18 totals, the SUMMARY block in G24/G36/I36) - recalculate
# correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DECEMBER 21")
$ws.Activate()

# Mark these tenants as paid in full (was blank -> equal to what was owed,
# zeroing out their balance in column I).
$ws.Range("H7").Value  = 4000
$ws.Range("H12").Value = 9000
$ws.Range("H13").Value = 8000

# Minimize the workbook window, and leave the cell cursor where the author
# left off reviewing the fix.
$excel.WindowState = -4140
$ws.Range("H14").Select()
